$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$gValues = @{
    2 = 0
    3 = 4
    4 = 2
    5 = 1
    6 = 3
    7 = 0
    8 = 1
    9 = 1
    10 = 2
    11 = 1
    12 = 0
    14 = 1
    15 = 2
    16 = 3
    17 = 1
    18 = 1
    19 = 1
    20 = 2
    21 = 1
    22 = 0
    23 = 1
    24 = 3
    25 = 0
    26 = 2
    27 = 2
    28 = 0
    29 = 1
    30 = 2
    31 = 1
    32 = 2
    33 = 1
    34 = 0
    35 = 1
    36 = 4
    37 = 1
    39 = 1
    40 = 1
    41 = 1
    42 = 2
    43 = 2
    44 = 2
    45 = 1
    46 = 3
    47 = 2
    48 = 1
    49 = 3
    50 = 2
    51 = 2
    52 = 0
    53 = 1
    54 = 2
    55 = 1
    56 = 0
    57 = 1
    58 = 0
    59 = 1
    60 = 1
    61 = 1
    62 = 2
    63 = 2
    64 = 0
    65 = 1
    66 = 1
    67 = 3
    68 = 2
    69 = 1
    71 = 1
    72 = 1
    73 = 2
    74 = 2
    75 = 1
    76 = 3
}
foreach ($row in $gValues.Keys) {
    $ws.Range("G$row").Value = $gValues[$row]
}
